$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update house_price_dollars values (column B) per linear regression v2 data
$ws.Range("B2").Value = 550
$ws.Range("B3").Value = 452
$ws.Range("B4").Value = 553
$ws.Range("B5").Value = 654
$ws.Range("B6").Value = 755
$ws.Range("B7").Value = 856
$ws.Range("B8").Value = 957
$ws.Range("B9").Value = 1058
$ws.Range("B10").Value = 1159
$ws.Range("B11").Value = 1060

# Update selection to just B2 (active cell only, no longer the full range)
$ws.Range("B2").Select()
